# Update the nowcasts worksheet with the results of the latest model run:
# refresh the numeric estimates for the existing rows and append a new
# row for the 2025-08-30 vintage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write columns B:K (10 numeric values) for a given row number.
function Set-DataRow($Row, $Values) {
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, $i + 2).Value = $Values[$i]
    }
}

# Helper: write a column-A label ("Row"/date) as literal text so Excel
# does not reinterpret date-like strings (e.g. "2025-08-30") as a date
# serial number. Uses a donor cell's (default) style so no stray
# "Text" number-format style is left behind on the written cell.
$donorStyle = $ws.Cells.Item(1, 2).Style
function Set-TextLabel($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $donorStyle
}

Set-DataRow 2 @([double]"0.35529062735081218", 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-DataRow 3 @([double]"0.33322932488502677", 0, [double]"-0.0051499207027774843", [double]"-0.00011291342515275761", [double]"-0.00093790221080953108", [double]"0.00031773024156445454", [double]"-6.9651418738962315e-05", [double]"-0.00072083949097101398", 0, [double]"0.007718400448858298")
Set-DataRow 4 @([double]"0.32150701032021306", [double]"-0.0023018916317374179", 0, [double]"0.00047772729980039458", [double]"6.868016734961851e-05", 0, [double]"0.00014557860882880879", [double]"-0.002587529728343248", [double]"0.00034046964700825186", [double]"-0.00076610555806710945")
Set-DataRow 5 @([double]"0.35580467549320083", [double]"0.0095099242951275625", [double]"-0.0074568509218240261", [double]"0.00011373266260338435", [double]"0.00099172917728496384", [double]"-0.0015810917159667689", [double]"1.889332165193007e-05", [double]"-0.0007083506336505287", 0, [double]"-0.00027220378905751241")
Set-DataRow 6 @([double]"0.37445452351912345", [double]"0.027539534987914022", 0, [double]"-0.00032496654339583255", [double]"1.9779760689569702e-05", 0, [double]"-6.3907821354194466e-05", [double]"-0.0022597581549458979", 0, [double]"0.00018595980891306096")
Set-DataRow 7 @([double]"0.29690104460190275", 0, [double]"-0.0026499943018211125", [double]"-0.00089064459055649564", [double]"-0.0062648470923306636", [double]"0.0010815367443699887", 0, [double]"0.00024170046862809013", 0, [double]"-0.0013590867449573829")
Set-DataRow 8 @([double]"0.16744957426771001", [double]"-0.058417806521557675", 0, [double]"-3.9245941296442858e-05", [double]"-0.00036194307868962273", 0, [double]"5.1309491888608729e-05", [double]"0.00089158522824544606", 0, [double]"-0.00035642031431160071")
Set-DataRow 9 @([double]"0.19590526220428095", 0, [double]"0.0044323024119971993", [double]"-0.0033204286031518928", [double]"-0.0045295828881494436", [double]"0.0013386679293828461", [double]"-0.00015324902492267335", [double]"-0.0001862869139106772", 0, [double]"-0.00053235815376229123")
Set-DataRow 10 @([double]"0.40061115606382292", [double]"0.099226985922044306", 0, [double]"-0.000365166072996234", [double]"-0.00033055372627747022", 0, [double]"-1.7648798637545545e-05", [double]"0.00035830173599825968", [double]"-0.0021980118939890054", [double]"-0.00041305182915371552")
Set-DataRow 11 @([double]"0.37784827511362706", 0, [double]"-0.027243915801356412", [double]"0.0024337296980999992", [double]"0.005798210556820717", [double]"0.0025414091335063771", [double]"0.00065072136068821328", [double]"0.0016258226457821284", 0, [double]"0.0057862865820084153")
Set-DataRow 12 @([double]"0.16134582824200724", [double]"-0.068106350395210971", 0, [double]"0.00036017348074736134", [double]"2.3928002410642268e-05", 0, [double]"1.2665756102842452e-05", [double]"-0.0019037968420538519", 0, [double]"-0.0042211796449230787")

# New vintage row (row 12) - add the date label and the full data row.
Set-TextLabel $ws.Cells.Item(12, 1) "2025-08-30"

# Column J got a bit narrower in this run.
$ws.Columns.Item(10).ColumnWidth = 14.333333333333334
